# 其他有價證券 (Other securities) sheet - sheet index 5 (1-based) - add the
# "otherbonds" record's full property_category/category/date/legislator_name/
# legislator_id/source_file/index columns (H:N), matching the template used by
# the 股票 (stock) sheet, and flesh out the header row to match too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# ---- Header row (row 1): B1:G1 change meaning, H1:N1 are brand new ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "face_value"
$ws.Range("F1").Value = "currency"
$ws.Range("G1").Value = "total"

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# New header cells need the same bold / bordered / centered look as the rest
# of row 1 (style index "1" in the original file).
$hdrNew = $ws.Range("H1:N1")
$hdrNew.Font.Bold = $true
$hdrNew.HorizontalAlignment = -4108
$hdrNew.VerticalAlignment = -4160
$hdrNew.Borders.LineStyle = 1

# ---- Data row (row 2) ----
# F2/G2 swap meaning: F2 becomes the currency string, G2 becomes the numeric total.
$ws.Range("F2").Value = "新臺幣"
$ws.Range("G2").Value = 10

# H2:N2 are brand new - same values/order as the 股票 sheet's matching row,
# except property_category = "otherbonds".
$ws.Range("H2").Value = "otherbonds"
$ws.Range("I2").Value = "normal"
# "2013-11-01" must stay literal text (like the source spreadsheet's shared
# string), not get auto-converted into a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2013-11-01"
$ws.Range("K2").Value = "鄭汝芬"
$ws.Range("L2").Value = 1713
$ws.Range("M2").Value = "tmp4cfc1"
$ws.Range("N2").Value = 187
